$wb = $excel.ActiveWorkbook

# --- Sheet "TP" -> "LL" ---
$wsLL = $wb.Worksheets.Item("TP")
$wsLL.Name = "LL"
$wsLL.Range("E1").Value = "LLRating"

# --- Sheet "FP" -> "NL" ---
$wsNL = $wb.Worksheets.Item("FP")
$wsNL.Name = "NL"
$wsNL.Range("D1").Value = "NLRating"

# --- Sheet "TRUTH": drop the Paradigm / FROC / FCTRL columns (D:F) ---
$wsTruth = $wb.Worksheets.Item("TRUTH")
$wsTruth.Columns("D:F").Delete()

# --- Restore the cursor/selection positions recorded in each sheet ---
$wsLL.Activate()
$wsLL.Range("G21").Select()

$wsNL.Activate()
$wsNL.Range("G15").Select()

$wsTruth.Activate()
$wsTruth.Range("G12").Select()
